$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-25 Monday", "2024-11-26 Tuesday"),
    @("54×18=", "93×99="),
    @("29×73=", "72×49="),
    @("62×42=", "75×26="),
    @("84×78=", "43×69="),
    @("12×21=", "57×41="),
    @("84×51=", "36×16="),
    @("80×51=", "28×27="),
    @("72×50=", "97×26="),
    @("26×51=", "81×97="),
    @("94×86=", "43×96="),
    @("60×37=", "50×59="),
    @("96×79=", "58×19="),
    @("14×30=", "59×77="),
    @("13×74=", "57×85="),
    @("74×53=", "38×25="),
    @("19×68=", "84×74="),
    @("74×88=", "18×76="),
    @("93×27=", "72×11="),
    @("78×44=", "92×95="),
    @("36×20=", "76×69="),
    @("89×80=", "62×22="),
    @("38×41=", "88×78="),
    @("88×27=", "41×14="),
    @("53×53=", "52×49="),
    @("31×38=", "67×71=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
